$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SHA")

# --- Build row 7 as a copy of row 6's formatting, with new content ---

# 1) Copy formatting (fill/border/alignment/number format) from row 6 to row 7
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Match row height to row 6 (234.75, custom height)
$ws.Rows.Item(7).RowHeight = $ws.Rows.Item(6).RowHeight

# 3) Fill in the new tool's data
$ws.Cells.Item(7,1).Value2 = "ClusterHiveReaderFromClusterLog.ps1"
$ws.Cells.Item(7,2).Value2 = @'
Purpose: 
# SYNOPSIS ClusterHiveReaderFromClusterLog.ps1
    The script Reads the Header Lines from a Server 2016 or later cluster.log and creates config files
    Script Name:  ClusterHiveReaderFromClusterLog.ps1    	
    Version:      1.1
    Last Update:  12 Feb 2020
    Author:       Josef Holzer 
## DESCRIPTION
    In 2016 and later we store all Configuration Info on top of the Cluster Log in csv format
    The script reads the header of the cluster.log file and creates the following files:    
    ClusterLogName-ClusConfig-All.txt             # Contains all Info Exported as AllObjects  fl *
    ClusterLogName-ClusConfig-All.xml             # Contains all Info as Powershell Objects
    ClusterLogName-ClusConfig-All-Overview.txt    # Contains most important info as | ft Prop1, Prop2...
    ClusterLogName-ClusConfig-ProcessIDs.txt      # All PIDs of Processes that wrote to cluster log
EXAMPLE 1
     ClusterHiveReaderFromClusterLog.ps1    
    - if you have several cluster logs in c:\logs and copy the script to this folder
       you simply run the script with no parameters
    - it will then take the first (2016 or later ) 
      cluster.log that contains config data and processes it
EXAMPLE 2
    ClusterHiveReaderFromClusterLog.ps1 -Path "C:\ClusterLog\H19N1.H19Corp.com_cluster.log" `
    -FileWithProcessInfoPathFull "C:\ClusterLog\H19N1-GeneralInfoPerHost.xml" -FindPIDs $True
    ...it will read cluster configuration and write it down into file names mentioned above including *ProcessIDs
'@
$ws.Cells.Item(7,5).Value2 = "https://github.com/CSS-Windows/WindowsDiag/blob/master/SHA/ClusterHiveReaderFromClusterLog/ClusterHiveReaderFromClusterLog.zip?raw=true"
$ws.Cells.Item(7,6).Value2 = $ws.Cells.Item(6,6).Value2

# 4) Add the hyperlink on E7 pointing at the zip download
$url = "https://github.com/CSS-Windows/WindowsDiag/blob/master/SHA/ClusterHiveReaderFromClusterLog/ClusterHiveReaderFromClusterLog.zip?raw=true"
$ws.Hyperlinks.Add($ws.Cells.Item(7,5), $url)

# Adding the hyperlink resets the cell style, so re-apply E6's formatting to E7
$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the view to reflect the new last row (mirrors the authored workbook) ---
$ws.Activate() | Out-Null
$ws.Range("E7").Select() | Out-Null
